# Add new "Network parameters" section (Batch Size, Exclude Unlabelled,
# Network parameters header, Pooling time ratio, RNN layers, RNN Cells,
# CNN layers, CNN filters, dropout) below the existing learning-parameter
# rows, mirroring the author's commit: "Added synth/unlabel. Remove
# dependency on unlabel dataset when excluded".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows under the original learning-rate block (rows 1,3-6 already exist)
$ws.Range("A7").Value  = "Batch Size"
$ws.Range("A8").Value  = "Exclude Unlabelled"

# New "Network parameters" section header
$ws.Range("A10").Value = "Network parameters"

# Network parameter rows (label in column A, value in column B)
$ws.Range("A12").Value = "Pooling time ratio"
$ws.Range("B12").Value = 4

$ws.Range("A13").Value = "RNN layers"
$ws.Range("B13").Value = 2

$ws.Range("A14").Value = "RNN Cells"
$ws.Range("B14").Value = 128

$ws.Range("A15").Value = "CNN layers"
$ws.Range("B15").Value = 6

$ws.Range("A16").Value = "CNN filters"
$ws.Range("B16").Value = "16,32,64,128,128,128"

$ws.Range("A17").Value = "dropout"
$ws.Range("B17").Value = 0.5

# Match column B's auto-fit width from the new longer values
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Move the active selection to D10, matching the saved workbook state
$ws.Range("D10").Select() | Out-Null
